$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Step 1: Rewrite the "Cada funcionalidade terá uma branch" paragraph into the
# new, longer paragraph text about initial commits going to master.
# ---------------------------------------------------------------------------
$newP4Text = "Os commits inicias (antes do início do desenvolvimento) serão todos feitos na master, a partir do momento da criação das primeiras branches, não deve ser feito mais nenhum commit diretamente na master."
$d.Content.Find.Execute("Cada funcionalidade terá uma branch", $true, $false, $false, $false, $false, $true, 1, $false, $newP4Text, 2) | Out-Null

# ---------------------------------------------------------------------------
# Step 2: Insert a brand-new list paragraph right after it:
# "Cada usuários do grupo terá uma branch" carrying the "_GoBack" bookmark at
# the very end of its text (this mirrors the bookmark that used to sit on the
# "Pagar uma Heineken..." paragraph - Word only ever keeps one "_GoBack", so
# re-adding it here automatically relocates it off of the old paragraph).
# ---------------------------------------------------------------------------
$r = $d.Content
$r.Find.ClearFormatting()
$r.Find.Execute($newP4Text) | Out-Null
$p4 = $r.Paragraphs(1)
$p4.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs(5)
$newPara.Range.Text = "Cada usuários do grupo terá uma branch"

# Work around the Range.Collapse() edge case that snaps a collapsed range back
# to the start of its paragraph whenever the target offset sits exactly at the
# paragraph's text end: temporarily append a placeholder character, anchor the
# bookmark immediately before it via Find, then remove the placeholder again.
$newParaRange = $d.Paragraphs(5).Range
$newParaRange.MoveEnd(1, -1) | Out-Null
$newParaRange.InsertAfter("X")

$bmRange = $d.Content
$bmRange.Find.ClearFormatting()
$bmRange.Find.Execute("Cada usuários do grupo terá uma branch") | Out-Null
$bmRange.Collapse(0)
$d.Bookmarks.Add("_GoBack", $bmRange)

$placeholder = $d.Range($bmRange.Start, $bmRange.Start + 1)
$placeholder.Delete()

Write-Output "Done. Paragraph count now: $($d.Paragraphs.Count)"
